# Auto-generated edit script: updates crypto price/volume table values
# to match the refreshed data pulled by the scheduled GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.931.88"
$ws.Range("E2").Value = "  -1.26%  "
$ws.Range("D3").Value = "2.353.59"
$ws.Range("E3").Value = "  -1.36%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'504.26"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").Value = "'130.15"
$ws.Range("E6").Value = "  -1.71%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  -2.53%  "
$ws.Range("D9").Value = "2.366.40"
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("D10").Value = "'0.0971"
$ws.Range("E10").Value = "  -0.35%  "
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("D12").Value = "'4.79"
$ws.Range("E12").Value = "  +2.55%  "
$ws.Range("E13").Value = "  -1.24%  "
$ws.Range("D14").Value = "2.771.96"
$ws.Range("E14").Value = "  -1.47%  "
$ws.Range("D15").Value = "55.912.20"
$ws.Range("E15").Value = "  -1.20%  "
$ws.Range("D16").Value = "'21.45"
$ws.Range("E16").Value = "  -1.04%  "
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("D18").Value = "2.390.53"
$ws.Range("E18").Value = "  +0.58%  "
$ws.Range("D19").Value = "'9.90"
$ws.Range("E19").Value = "  -2.79%  "
$ws.Range("D20").Value = "'310.79"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("E21").Value = "  -0.89%  "
$ws.Range("E22").Value = "  -1.05%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").Value = "'65.43"
$ws.Range("E24").Value = "  -1.24%  "
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("D27").Value = "'0.145"
$ws.Range("E27").Value = "  -2.92%  "
$ws.Range("D28").Value = "'7.13"
$ws.Range("E28").Value = "  -3.02%  "
$ws.Range("D29").Value = "'170.96"
$ws.Range("E29").Value = "  -2.66%  "
$ws.Range("D30").Value = "0.0₃0704"
$ws.Range("E30").Value = "  -2.84%  "
$ws.Range("D31").Value = "'1.63"
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "'0.996"
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D34").Value = "'5.73"
$ws.Range("E34").Value = "  -2.41%  "
$ws.Range("D35").Value = "'1.06"
$ws.Range("E35").Value = "  -4.77%  "
$ws.Range("D36").Value = "'17.66"
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("E37").Value = "  -2.07%  "
$ws.Range("D38").Value = "'0.836"
$ws.Range("E38").Value = "  +2.19%  "
$ws.Range("E39").Value = "  -4.02%  "
$ws.Range("D40").Value = "'36.18"
$ws.Range("E40").Value = "  -1.66%  "
$ws.Range("E41").Value = "  -3.01%  "
$ws.Range("D42").Value = "'3.34"
$ws.Range("E42").Value = "  -1.08%  "
$ws.Range("D43").Value = "'4.85"
$ws.Range("E43").Value = "  +0.66%  "
$ws.Range("D44").Value = "'125.71"
$ws.Range("E44").Value = "  -4.92%  "
$ws.Range("D45").Value = "'0.557"
$ws.Range("E45").Value = "  -1.74%  "
$ws.Range("E46").Value = "  -1.09%  "
$ws.Range("D47").Value = "'240.60"
$ws.Range("E47").Value = "  -2.33%  "
$ws.Range("D48").Value = "'0.0477"
$ws.Range("E48").Value = "  -1.34%  "
$ws.Range("E49").Value = "  -1.23%  "
$ws.Range("D50").Value = "'0.0206"
$ws.Range("E50").Value = "  -1.77%  "
$ws.Range("D51").Value = "'16.57"
$ws.Range("E51").Value = "  -3.43%  "
